# Refactor authentication flow: update login, logout, and session management;
# enhance login.jsp with error handling and styling; add dashboard.jsp;
# implement CSS for improved UI.
#
# Spreadsheet-side change: log a new progress entry (row 10) on the
# "Tiến trình" (progress) sheet, recolor the previous entry's priority
# marker now that the login bug is fixed, and move the viewport/selection
# down to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Update row 9's priority marker: the login session bug it flagged is now
#     fixed, so recolor its priority cell from red ("cao") to yellow.
$ws.Range("E9").Interior.Color = 65535

# --- Add the new progress row (row 10), copying row 9's look & feel for the
#     date/content columns (thin border, same number format/wrap).
$ws.Range("A9:D9").Copy($ws.Range("A10:D10"))
$ws.Rows(10).RowHeight = 43.2

$ws.Range("A10").Value = 46008
$ws.Range("B10").Value = "Đaã sửa lỗi đăng nhập, kiểm tra đăng nhập hoàn tất"
$ws.Range("C10").Value = "thực hiện đăng ký và xoá người dùng"
$ws.Range("D10").Value = "Trạng thái đăng nhập khi sai mật khẩu vẫn trả về user để so sánh key nên sai, đã sửa lại ở mục UserService"

# Priority cell for the new row: red, but with no border underneath (last row).
$e10 = $ws.Range("E10")
$e10.Interior.Color = 192
$e10.Borders.LineStyle = 0

# --- Move the on-screen selection down to where the new row sits.
$ws.Range("H12").Select()
